# Update the "想去人数" (wanted-attendance) counts on the 展览 (Exhibition)
# and 全部类型 (All Types) sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 824
$wsExhibition.Range("F4").Value = 284
$wsExhibition.Range("F5").Value = 1014
$wsExhibition.Range("F6").Value = 2376

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 824
$wsAll.Range("F4").Value = 284
$wsAll.Range("F7").Value = 1014
$wsAll.Range("F8").Value = 2376
